# Create a new "TeamInfo" worksheet (after "SoccerPage") holding the Real
# Madrid assertion-category labels used by the new test cases, and make it
# the active/selected sheet (mirrors the "created test cases to retrieve
# stats and team info" commit).

$wb = $excel.ActiveWorkbook

# Duplicate the existing sheet so the new one inherits the same look/format
# (column formatting, view settings, namespaces, etc.), then rename it and
# wipe its content before filling in the new data. Copy() places the
# duplicate right after SoccerPage and leaves it as the active sheet, so
# grab it from there instead of guessing its auto-generated name.
$soccerPage = $wb.Worksheets.Item("SoccerPage")
$soccerPage.Copy($null, $soccerPage)
$teamInfo = $wb.ActiveSheet
$teamInfo.Name = "TeamInfo"
$teamInfo.Cells.Clear()

$values = @(
    "Assertions",
    "Real Madrid Performance Stats",
    "Real Madrid Discipline Stats",
    "Real Madrid Scoring Stats",
    "Real Madrid Results",
    "Real Madrid Squad",
    "Real Madrid Transfers"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $teamInfo.Cells.Item($i + 1, 1).Value = $values[$i]
}

$teamInfo.Columns.Item(1).AutoFit()

# Make the new sheet the active tab/selection, like it was left after the edit.
$teamInfo.Select()
$teamInfo.Range("A8").Select()
